# "Them ham xep theo alpha, chinh tong alpha = 12"
# - Rename "Web" -> "WEB" in two course names.
# - Recompute the "alpha" totals: numbers that used to be stored as
#   numeric cells are replaced by the *text* rendering of the new total
#   (the sheet now stores these tallies as text, e.g. via a helper
#   "xep theo alpha" formula whose result got pasted as text).
# - A couple of brand-new totals appear (E6, F16).
# - B2 (previously blank) now shows the text "12".
# - Column A is widened, and the live selection moves to J8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the two "Web" -> "WEB" rows, keeping their existing
#     (quote-prefixed) cell style untouched. Prefixing the literal with
#     an apostrophe lets the host know this is text without us having to
#     touch NumberFormat/Style (these two cells already carry the
#     workbook's "text" style).
$ws.Range("A4").Value = "'Công nghệ thiết kế WEB"
$ws.Range("A5").Value = "'Công nghệ thiết kế WEB nâng cao"

# --- Numeric "alpha" tallies that become text cells (same displayed
#     digits except the totals that grew from 5/3/10 -> 6/4/12), plus two
#     brand-new text cells (E6, F16). Each value is written with a
#     leading apostrophe so the host stores it as literal text instead of
#     re-parsing the digits back into a number, then the cell style is
#     reset to "Normal" so we don't leave a stray quote-prefix style
#     behind (these cells were plain General-formatted before the edit).
$alphaCells = @(
    @("B2","12"),
    @("B3","6"),
    @("E3","6"),
    @("B4","4"),
    @("E4","4"),
    @("F4","4"),
    @("B5","6"),
    @("F5","6"),
    @("B6","6"),
    @("E6","6"),
    @("B7","6"),
    @("B8","6"),
    @("C9","12"),
    @("C10","6"),
    @("F10","6"),
    @("C11","6"),
    @("D11","6"),
    @("C12","12"),
    @("D13","12"),
    @("G14","12"),
    @("F15","12"),
    @("F16","6"),
    @("G16","6"),
    @("H16","6"),
    @("G17","4"),
    @("H17","4"),
    @("G18","6"),
    @("H18","6"),
    @("G19","4"),
    @("H19","4"),
    @("H20","12"),
    @("I21","12"),
    @("I22","12"),
    @("F23","12")
)

foreach ($pair in $alphaCells) {
    $addr = $pair[0]
    $text = $pair[1]
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# --- Widen column A (new "xep theo alpha" column gets more room).
$ws.Columns.Item(1).ColumnWidth = 51.1666666666667

# --- Move the live selection to J8.
$ws.Range("J8").Select() | Out-Null
